# "new Madigan bike hours" - shift the Driver Position 5 timepoints schedule
# forward onto the new bike-hours timetable and drop the trailing trips that
# no longer exist (rows 72-91).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trips that only change their Departure time (rows 4-27)
$ws.Range("D4:D7").Value   = "18:16:39"
$ws.Range("D8:D11").Value  = "18:33:18"
$ws.Range("D12:D15").Value = "18:48:18"
$ws.Range("D16:D19").Value = "19:13:18"
$ws.Range("D20:D23").Value = "19:38:18"
$ws.Range("D24:D27").Value = "19:54:57"

# Trips that shift to a different stop in the cycle as well as a new time (rows 28-71)
$ws.Range("C28:C31").Value = "to Lewis North"
$ws.Range("D28:D31").Value = "20:09:57"

$ws.Range("A32:A35").Value = "480a"
$ws.Range("B32:B35").Value = "Warrior Zone"
$ws.Range("D32:D35").Value = "20:34:57"

$ws.Range("C36:C39").Value = "to Passenger Terminal"
$ws.Range("D36:D39").Value = "20:59:57"

$ws.Range("A40:A43").Value = "772b"
$ws.Range("B40:B43").Value = "Evergreen McChord Lodging"
$ws.Range("D40:D43").Value = "21:16:36"

$ws.Range("C44:C47").Value = "to Lewis North"
$ws.Range("D44:D47").Value = "21:31:36"

$ws.Range("A48:A51").Value = "480a"
$ws.Range("B48:B51").Value = "Warrior Zone"
$ws.Range("D48:D51").Value = "21:56:36"

$ws.Range("C52:C55").Value = "to Passenger Terminal"
$ws.Range("D52:D55").Value = "22:21:36"

$ws.Range("A56:A59").Value = "772b"
$ws.Range("B56:B59").Value = "Evergreen McChord Lodging"
$ws.Range("D56:D59").Value = "22:38:15"

$ws.Range("C60:C63").Value = "to Lewis North"
$ws.Range("D60:D63").Value = "22:53:15"

$ws.Range("A64:A67").Value = "480a"
$ws.Range("B64:B67").Value = "Warrior Zone"
$ws.Range("D64:D67").Value = "23:18:15"

$ws.Range("C68:C71").Value = "to Passenger Terminal"
$ws.Range("D68:D71").Value = "23:43:15"

# The schedule is now only 17 trips (was 22) - remove the trailing 5 trips / 20 rows
$ws.Rows("72:91").Delete()
